# Updating writing excel code
#
# 1. Drop the stale "Backup" sheet (it was just a duplicate copy of
#    "Search Items").
# 2. Add a new "TestSheet" worksheet (after "Search Items") with some
#    Name/City/Status test data.
# 3. Add a "Status" column with the QA result ("PASSED") to "Search Items".

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# --- 1. remove the old "Backup" sheet -------------------------------------
$backup = $wb.Worksheets.Item("Backup")
$backup.Delete()

$searchItems = $wb.Worksheets.Item("Search Items")

# --- 2. add the new "TestSheet" worksheet and fill it in -------------------
$testSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $searchItems)
$testSheet.Name = "TestSheet"

$testSheet.Range("A1").Value = "Name"
$testSheet.Range("B1").Value = "City"
$testSheet.Range("C1").Value = "Status"

$testSheet.Range("A2").Value = "Richa"
$testSheet.Range("B2").Value = "Chandigarh"

$testSheet.Range("A3").Value = "Shikha"
$testSheet.Range("B3").Value = "Gurgaon"

$testSheet.Range("C2").Value = "Passed"
$testSheet.Range("C3").Value = "Passed"

# copy formatting for the header row / data rows from "Search Items"
# (TestSheet's header uses the same plain bordered style as the data rows,
# not the yellow-highlighted header style)
$searchItems.Range("A2:B2").Copy()
$testSheet.Range("A1:C1").PasteSpecial(-4122)   # xlPasteFormats
$searchItems.Range("A2:B3").Copy()
$testSheet.Range("A2:B3").PasteSpecial(-4122)   # xlPasteFormats

$testSheet.Columns.Item(2).ColumnWidth = 13.88
$testSheet.Columns.Item(3).ColumnWidth = 14.45

# --- 3. extend "Search Items" with a Status column --------------------------
$searchItems.Range("C1").Value = "Status"
$searchItems.Range("A1").Copy()
$searchItems.Range("C1").PasteSpecial(-4122)   # xlPasteFormats

$searchItems.Range("C2").Value = "PASSED"
$searchItems.Range("C3").Value = "PASSED"

$searchItems.Columns.Item(3).ColumnWidth = 18.0

# make "Search Items" the active sheet/tab again
$searchItems.Select()
$searchItems.Range("I17").Select()
$testSheet.Range("C3").Select()
$searchItems.Select()
